# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking text (e.g. "0.999", "562.59").
# Force it to Text format before writing so Excel keeps these as strings
# instead of silently converting them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.536.83"
$ws.Range("E2").Value = "  -1.96%  "

$ws.Range("D3").Value = "2.530.87"
$ws.Range("E3").Value = "  -3.23%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "562.59"
$ws.Range("E5").Value = "  -2.22%  "

$ws.Range("D6").Value = "149.99"
$ws.Range("E6").Value = "  -3.85%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").Value = "0.610"
$ws.Range("E8").Value = "  -2.03%  "

$ws.Range("D9").Value = "2.528.58"
$ws.Range("E9").Value = "  -3.14%  "

$ws.Range("D10").Value = "0.111"
$ws.Range("E10").Value = "  -6.66%  "

$ws.Range("D11").Value = "5.51"
$ws.Range("E11").Value = "  -5.40%  "

$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("D13").Value = "0.368"
$ws.Range("E13").Value = "  -3.68%  "

$ws.Range("D14").Value = "26.99"
$ws.Range("E14").Value = "  -4.40%  "

$ws.Range("D15").Value = "2.973.02"
$ws.Range("E15").Value = "  -3.74%  "

$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  -6.36%  "

$ws.Range("D17").Value = "62.345.97"
$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("D18").Value = "2.498.85"
$ws.Range("E18").Value = "  -3.60%  "

$ws.Range("D19").Value = "11.46"
$ws.Range("E19").Value = "  -5.04%  "

$ws.Range("D20").Value = "7.22"
$ws.Range("E20").Value = "  -5.60%  "

$ws.Range("D21").Value = "4.33"
$ws.Range("E21").Value = "  -5.10%  "

$ws.Range("D22").Value = "328.81"
$ws.Range("E22").Value = "  -4.22%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").Value = "65.31"
$ws.Range("E24").Value = "  -3.09%  "

$ws.Range("D25").Value = "1.83"
$ws.Range("E25").Value = "  +3.32%  "

$ws.Range("D26").Value = "0.0000107"
$ws.Range("E26").Value = "  -2.01%  "

$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").Value = "1.59"
$ws.Range("E27").Value = "  +1.04%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.648.82"
$ws.Range("E28").Value = "  -3.58%  "

$ws.Range("D29").Value = "8.73"
$ws.Range("E29").Value = "  -5.13%  "

$ws.Range("D30").Value = "555.71"
$ws.Range("E30").Value = "  -6.10%  "

$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "7.98"
$ws.Range("E31").Value = "  +0.72%  "

$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.23%  "

$ws.Range("D33").Value = "0.155"
$ws.Range("E33").Value = "  -3.30%  "

$ws.Range("D34").Value = "1.95"
$ws.Range("E34").Value = "  -5.65%  "

$ws.Range("D35").Value = "1.64"
$ws.Range("E35").Value = "  -6.63%  "

$ws.Range("D36").Value = "6.11"
$ws.Range("E36").Value = "  -7.04%  "

$ws.Range("D37").Value = "5.02"
$ws.Range("E37").Value = "  -6.23%  "

$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("D39").Value = "0.389"
$ws.Range("E39").Value = "  -4.48%  "

$ws.Range("D40").Value = "18.93"
$ws.Range("E40").Value = "  -4.10%  "

$ws.Range("D41").Value = "150.77"
$ws.Range("E41").Value = "  -2.41%  "

$ws.Range("D42").Value = "1.77"
$ws.Range("E42").Value = "  -5.29%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "40.94"
$ws.Range("E44").Value = "  -1.33%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.38"
$ws.Range("E45").Value = "  -1.60%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "152.50"
$ws.Range("E46").Value = "  -2.08%  "

$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "3.74"
$ws.Range("E47").Value = "  -4.63%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "22.50"
$ws.Range("E48").Value = "  -3.73%  "

$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0557"
$ws.Range("E49").Value = "  -5.76%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.603"
$ws.Range("E50").Value = "  -3.99%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.0964"
$ws.Range("E51").Value = "  -5.13%  "
